$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the cells we touch so numeric-looking strings
# (e.g. "0.620", thousand-dotted prices) keep their exact original formatting
# instead of being auto-converted to numbers (which would drop trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.991.70"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.724.07"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.99%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "530.49"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.87"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.35%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.745.91"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.24"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +15.67%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.343"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.93%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.192.80"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.977.72"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.10%  "

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.59"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.45%  "

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.818.95"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +7.53%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "347.85"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.77%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.64"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.47"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.81%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.54"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.52%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +5.79%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.422"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.989"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0831"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.37"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.81"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +10.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.19"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.02"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.29"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +8.41%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +10.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.929"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.907"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +9.21%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +9.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.28"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.71"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "285.32"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.54"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.06%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.620"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.97%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0993"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.36%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.135.47"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +8.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.99"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +7.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0547"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.22%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.54"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +7.04%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.53"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.18%  "
